# chore: update Sheets via scheduled runner
# Refreshes cached market-board price/profit figures (columns H-N) across
# the per-job leve sheets with the latest pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1664.762
$ws.Range("J17").Value = 1664.762
$ws.Range("L17").Value = 4994.286
$ws.Range("N17").Value = -5330.286

$ws.Range("H70").Value = 5021.769
$ws.Range("J70").Value = 5831.8887
$ws.Range("L70").Value = 17495.6661
$ws.Range("N70").Value = -18035.6661

$ws.Range("H73").Value = 5021.769
$ws.Range("J73").Value = 5831.8887
$ws.Range("L73").Value = 17495.6661
$ws.Range("N73").Value = -19367.6661

$ws.Range("H88").Value = 7615.385
$ws.Range("J88").Value = 7749.75
$ws.Range("L88").Value = 7749.75
$ws.Range("N88").Value = -8561.75

$ws.Range("H91").Value = 7615.385
$ws.Range("J91").Value = 7749.75
$ws.Range("L91").Value = 7749.75
$ws.Range("N91").Value = -10557.75

$ws.Range("H98").Value = 1017.7857
$ws.Range("I98").Value = 1019.1539
$ws.Range("K98").Value = 1019.1539
$ws.Range("M98").Value = 478.8461

$ws.Range("H113").Value = 3474.7
$ws.Range("I113").Value = 3124.5
$ws.Range("K113").Value = 3124.5
$ws.Range("M113").Value = 129.5

$ws.Range("H116").Value = 38473572
$ws.Range("I116").Value = 55568056
$ws.Range("J116").Value = 10986.75
$ws.Range("K116").Value = 55568056
$ws.Range("L116").Value = 10986.75
$ws.Range("M116").Value = -55564614
$ws.Range("N116").Value = -17870.75

$ws.Range("H122").Value = 1017.7857
$ws.Range("I122").Value = 1019.1539
$ws.Range("K122").Value = 3057.4617
$ws.Range("M122").Value = -607.4616999999998

$ws.Range("H132").Value = 110174.53
$ws.Range("I132").Value = 129679.82
$ws.Range("K132").Value = 389039.46
$ws.Range("M132").Value = -386509.46

$ws.Range("H137").Value = 2708.5789
$ws.Range("I137").Value = 2435.0425
$ws.Range("K137").Value = 7305.127500000001
$ws.Range("M137").Value = -4755.127500000001

$ws.Range("H141").Value = 2239.8923
$ws.Range("J141").Value = 4358.8335
$ws.Range("L141").Value = 13076.5005
$ws.Range("N141").Value = -23436.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15565.91
$ws.Range("I32").Value = 18143.49
$ws.Range("K32").Value = 18143.49
$ws.Range("M32").Value = -17856.49

$ws.Range("H41").Value = 18559
$ws.Range("I41").Value = 5451.6665
$ws.Range("J41").Value = 31666.334
$ws.Range("K41").Value = 5451.6665
$ws.Range("L41").Value = 31666.334
$ws.Range("M41").Value = -5037.6665
$ws.Range("N41").Value = -32494.334

$ws.Range("H61").Value = 260690.11
$ws.Range("I61").Value = 3584.5
$ws.Range("K61").Value = 3584.5
$ws.Range("M61").Value = -3372.5

$ws.Range("H74").Value = 4263648.5

$ws.Range("H77").Value = 4263648.5

$ws.Range("H132").Value = 14519.833
$ws.Range("I132").Value = 19196.875
$ws.Range("K132").Value = 57590.625
$ws.Range("M132").Value = -55060.625

$ws.Range("H136").Value = 260690.11
$ws.Range("I136").Value = 3584.5
$ws.Range("K136").Value = 10753.5
$ws.Range("M136").Value = -8203.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1279.4166
$ws.Range("I99").Value = 1026.579
$ws.Range("J99").Value = 2240.2
$ws.Range("K99").Value = 1026.579
$ws.Range("L99").Value = 2240.2
$ws.Range("M99").Value = 471.421
$ws.Range("N99").Value = -5236.2

$ws.Range("H128").Value = 6487.5
$ws.Range("I128").Value = 6487.5
$ws.Range("K128").Value = 19462.5
$ws.Range("M128").Value = -16972.5

$ws.Range("H134").Value = 2217.75
$ws.Range("I134").Value = 1277.2821
$ws.Range("K134").Value = 3831.8463
$ws.Range("M134").Value = -1296.8463

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21280202
$ws.Range("I31").Value = 37039856
$ws.Range("J31").Value = 4671.5
$ws.Range("K31").Value = 37039856
$ws.Range("L31").Value = 4671.5
$ws.Range("M31").Value = -37039561
$ws.Range("N31").Value = -5261.5

$ws.Range("H34").Value = 21280202
$ws.Range("I34").Value = 37039856
$ws.Range("J34").Value = 4671.5
$ws.Range("K34").Value = 37039856
$ws.Range("L34").Value = 4671.5
$ws.Range("M34").Value = -37039654
$ws.Range("N34").Value = -5075.5

$ws.Range("H99").Value = 15219.909
$ws.Range("J99").Value = 13166.333
$ws.Range("L99").Value = 13166.333
$ws.Range("N99").Value = -16162.333

$ws.Range("H126").Value = 15219.909
$ws.Range("J126").Value = 13166.333
$ws.Range("L126").Value = 39498.999
$ws.Range("N126").Value = -44438.999

$ws.Range("H132").Value = 31019416
$ws.Range("I132").Value = 34195508
$ws.Range("K132").Value = 102586524
$ws.Range("M132").Value = -102583994

$ws.Range("H134").Value = 942.8905999999999
$ws.Range("I134").Value = 923.14514
$ws.Range("K134").Value = 2769.43542
$ws.Range("M134").Value = -234.4354199999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 62598110
$ws.Range("J37").Value = 62598110
$ws.Range("L37").Value = 187794330
$ws.Range("N37").Value = -187794554

$ws.Range("H92").Value = 612.5
$ws.Range("I92").Value = 750
$ws.Range("J92").Value = 475
$ws.Range("K92").Value = 2250
$ws.Range("L92").Value = 1425
$ws.Range("M92").Value = -1002
$ws.Range("N92").Value = -3921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H113").Value = 1874.5
$ws.Range("I113").Value = 1999.4
$ws.Range("K113").Value = 1999.4
$ws.Range("M113").Value = 170.5999999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 49969
$ws.Range("J70").Value = 49969
$ws.Range("L70").Value = 49969
$ws.Range("N70").Value = -50509

$ws.Range("H73").Value = 49969
$ws.Range("J73").Value = 49969
$ws.Range("L73").Value = 49969
$ws.Range("N73").Value = -51841

$ws.Range("H122").Value = 12662.417
$ws.Range("I122").Value = 4500
$ws.Range("J122").Value = 14294.9
$ws.Range("K122").Value = 13500
$ws.Range("L122").Value = 42884.7
$ws.Range("M122").Value = -11050
$ws.Range("N122").Value = -47784.7

$ws.Range("H132").Value = 1868.14
$ws.Range("I132").Value = 1886.8469
$ws.Range("J132").Value = 951.5
$ws.Range("K132").Value = 5660.5407
$ws.Range("L132").Value = 2854.5
$ws.Range("M132").Value = -3130.5407
$ws.Range("N132").Value = -7914.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2166.3333
$ws.Range("I62").Value = 1750
$ws.Range("J62").Value = 2999
$ws.Range("K62").Value = 1750
$ws.Range("L62").Value = 2999
$ws.Range("M62").Value = -1126
$ws.Range("N62").Value = -4247

$ws.Range("H65").Value = 2166.3333
$ws.Range("I65").Value = 1750
$ws.Range("J65").Value = 2999
$ws.Range("K65").Value = 8750
$ws.Range("L65").Value = 14995
$ws.Range("M65").Value = -5630
$ws.Range("N65").Value = -21235
